$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.409494638442993
$ws.Range("B1").Value = 2.75522518157959
$ws.Range("C1").Value = 2.920884609222412
$ws.Range("D1").Value = 3.556612253189087
$ws.Range("E1").Value = 0.8769075870513916
